$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.697.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "'1.896.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.12%  "
$ws.Range("D5").Value = "'312.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("E6").Value = "  -1.00%  "
$ws.Range("D7").Value = "'0.4870"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.86%  "
$ws.Range("E8").Value = "  -0.74%  "
$ws.Range("D9").Value = "'0.07324"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("D10").Value = "'0.9141"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'20.56"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.29%  "
$ws.Range("D12").Value = "'0.07659"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("D13").Value = "'1.893.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("D14").Value = "'5.482"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("D15").Value = "'6.613"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "'91.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("D18").Value = "'0.000008776"
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("D20").Value = "'27.741.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("E21").Value = "  -2.56%  "
$ws.Range("D22").Value = "'5.119"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "'2.113.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").Value = "'10.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("E25").Value = "  -1.95%  "
$ws.Range("D26").Value = "'153.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.61%  "
$ws.Range("E27").Value = "  -1.05%  "
$ws.Range("D28").Value = "'2.157"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.11%  "
$ws.Range("D29").Value = "'115.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("D30").Value = "'4.873"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.19%  "
$ws.Range("D31").Value = "'0.08897"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").Value = "'3.200"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.05%  "
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("D34").Value = "'0.7660"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.73%  "
$ws.Range("D35").Value = "'4.630"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.61%  "
$ws.Range("D36").Value = "'0.02039"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.26%  "
$ws.Range("D37").Value = "'2.530"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.39%  "
$ws.Range("D38").Value = "'1.093"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.19%  "
$ws.Range("D39").Value = "'0.05272"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.06%  "
$ws.Range("E40").Value = "  -2.55%  "
$ws.Range("D41").Value = "'2.979"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("D42").Value = "'6.890"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.36%  "
$ws.Range("D43").Value = "'8.523"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.36%  "
$ws.Range("E44").Value = "  -0.62%  "
$ws.Range("D45").Value = "'111.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.05%  "
$ws.Range("D46").Value = "'10.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("D47").Value = "'0.4784"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.96%  "
$ws.Range("E48").Value = "  -1.11%  "
$ws.Range("D49").Value = "'1.635"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.69%  "
$ws.Range("D50").Value = "'67.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("D51").Value = "'0.06051"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.18%  "
